# Uppercase the smaller-font "tail" portions of the two big section
# headings ("Unsupervised Learning" and "Dimensionality Reduction"),
# leaving the large drop-cap letters and everything else untouched.
#
# Each heading paragraph is split into runs like:
#   "U" (big) + "nsupervised" (small) + " L" (big) + "earning" (small)
# so the Find/Replace is scoped to each heading paragraph's own Range
# to avoid touching unrelated occurrences of "learning" etc. elsewhere
# in the document.

$d = $word.ActiveDocument

# --- "Unsupervised Learning" heading ---
$p1 = $d.Paragraphs(1).Range
$p1.Find.Execute("nsupervised", $true, $false, $false, $false, $false, $true, 1, $false, "NSUPERVISED", 2)

$p1 = $d.Paragraphs(1).Range
$p1.Find.Execute("earning", $true, $false, $false, $false, $false, $true, 1, $false, "EARNING", 2)

# --- "Dimensionality Reduction" heading ---
$p2 = $d.Paragraphs(33).Range
$p2.Find.Execute("imensionality", $true, $false, $false, $false, $false, $true, 1, $false, "IMENSIONALITY", 2)

$p2 = $d.Paragraphs(33).Range
$p2.Find.Execute("eduction", $true, $false, $false, $false, $false, $true, 1, $false, "EDUCTION", 2)
